$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "FT232500DKJ9M3LJ"
$ws.Range("A7").Value = "FT232500DR4M2YS7"
$ws.Range("A8").Value = "FT232500DXZ685YD"
